$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 0.006728519194957983
$ws.Range("D2").Value = 0.0008797679034433514
$ws.Range("E2").Value = 0.4328853915826016
$ws.Range("F2").Value = 0.5140666435330417
$ws.Range("G2").Value = 0.4148371388643
$ws.Range("H2").Value = 0.415906419037583
$ws.Range("I2").Value = 0.378528808120933
$ws.Range("O2").Value = 1.604195653508384
$ws.Range("C3").Value = 0.006076801449463431
$ws.Range("D3").Value = 0.0008078252901420058
$ws.Range("E3").Value = 0.3773509119814804
$ws.Range("F3").Value = 0.477166177030341
$ws.Range("G3").Value = 0.3751547748629775
$ws.Range("H3").Value = 0.401380075461347
$ws.Range("I3").Value = 0.3510449145856356
$ws.Range("O3").Value = 1.489284539608462
$ws.Range("C4").Value = 0.005679248823945215
$ws.Range("D4").Value = 0.000763964996567168
$ws.Range("E4").Value = 0.3433622218953474
$ws.Range("F4").Value = 0.4548086757983754
$ws.Range("G4").Value = 0.3509929734859725
$ws.Range("H4").Value = 0.3927130817850184
$ws.Range("I4").Value = 0.334390143325102
$ws.Range("O4").Value = 1.41967317052405
$ws.Range("C5").Value = 0.005517902011231257
$ws.Range("D5").Value = 0.0007461718022279129
$ws.Range("E5").Value = 0.3295368069885853
$ws.Range("F5").Value = 0.4457726492480845
$ws.Range("G5").Value = 0.3411973456554591
$ws.Range("H5").Value = 0.3892444681041241
$ws.Range("I5").Value = 0.3276582463111168
$ws.Range("O5").Value = 1.391541781391453
$ws.Range("C6").Value = 0.005491150427758384
$ws.Range("D6").Value = 0.0007432221500911496
$ws.Range("E6").Value = 0.327242564663635
$ws.Range("F6").Value = 0.444276729850813
$ws.Range("G6").Value = 0.3395738211245032
$ws.Range("H6").Value = 0.3886723230525035
$ws.Range("I6").Value = 0.3265437343539261
$ws.Range("O6").Value = 1.386884778666513
$ws.Range("C7").Value = 0.005677070170630572
$ws.Range("D7").Value = 0.0007637247050613638
$ws.Range("E7").Value = 0.3431756683244203
$ws.Range("F7").Value = 0.454686510413751
$ws.Range("G7").Value = 0.3508606625685928
$ws.Range("H7").Value = 0.3926660469550143
$ws.Range("I7").Value = 0.3342991321805897
$ws.Range("O7").Value = 1.419292828240373
$ws.Range("C8").Value = 0.006503268574064691
$ws.Range("D8").Value = 0.0008548980653255711
$ws.Range("E8").Value = 0.4137129092959952
$ws.Range("F8").Value = 0.5012808901624766
$ws.Range("G8").Value = 0.4011120680816163
$ws.Range("H8").Value = 0.4108452740905477
$ws.Range("I8").Value = 0.3690063707470514
$ws.Range("O8").Value = 1.564377203308595
$ws.Range("C9").Value = 0.00814401422202593
$ws.Range("D9").Value = 0.001036116788627695
$ws.Range("E9").Value = 0.5530185435665942
$ws.Range("F9").Value = 0.595056017112114
$ws.Range("G9").Value = 0.5013017407642621
$ws.Range("H9").Value = 0.4485062746539938
$ws.Range("I9").Value = 0.4388371859021873
$ws.Range("O9").Value = 1.856474353692022
$ws.Range("C10").Value = 0.009362015047557293
$ws.Range("D10").Value = 0.001170690267082719
$ws.Range("E10").Value = 0.6561296499827876
$ws.Range("F10").Value = 0.6654628728308296
$ws.Range("G10").Value = 0.5759689139247541
$ws.Range("H10").Value = 0.4774197779458689
$ws.Range("I10").Value = 0.4912553836842619
$ws.Range("O10").Value = 2.075855294684573
$ws.Range("C11").Value = 0.009918848095985311
$ws.Range("D11").Value = 0.001232215773638856
$ws.Range("E11").Value = 0.7032396335836353
$ws.Range("F11").Value = 0.6978308986636961
$ws.Range("G11").Value = 0.6101787504348692
$ws.Range("H11").Value = 0.4908474577540858
$ws.Range("I11").Value = 0.5153513634726465
$ws.Range("O11").Value = 2.17672900672693
$ws.Range("C12").Value = 0.01013010071406484
$ws.Range("D12").Value = 0.001255557420938302
$ws.Range("E12").Value = 0.7211111176288938
$ws.Range("F12").Value = 0.7101373296406308
$ws.Range("G12").Value = 0.6231689435713008
$ws.Range("H12").Value = 0.4959719445542419
$ws.Range("I12").Value = 0.5245124474168676
$ws.Range("O12").Value = 2.21508431669082
$ws.Range("C13").Value = 0.01008458632070841
$ws.Range("D13").Value = 0.001250528467938494
$ws.Range("E13").Value = 0.7172607109845899
$ws.Range("F13").Value = 0.7074847156648474
$ws.Range("G13").Value = 0.6203696757861223
$ws.Range("H13").Value = 0.4948665249933981
$ws.Range("I13").Value = 0.5225378159277341
$ws.Range("O13").Value = 2.206816818587299
$ws.Range("C14").Value = 0.009936220144879826
$ws.Range("D14").Value = 0.001234135242992096
$ws.Range("E14").Value = 0.7047092759485309
$ws.Range("F14").Value = 0.6988423638079126
$ws.Range("G14").Value = 0.6112467426643207
$ws.Range("H14").Value = 0.4912682552297838
$ws.Range("I14").Value = 0.5161043182866933
$ws.Range("O14").Value = 2.179881371418446
$ws.Range("C15").Value = 0.009845392495471117
$ws.Range("D15").Value = 0.001224099528862155
$ws.Range("E15").Value = 0.6970253985258097
$ws.Range("F15").Value = 0.693555117943859
$ws.Range("G15").Value = 0.605663348196174
$ws.Range("H15").Value = 0.4890693893524372
$ws.Range("I15").Value = 0.5121683756061799
$ws.Range("O15").Value = 2.163403085112861
$ws.Range("C16").Value = 0.009325679808874554
$ws.Range("D16").Value = 0.001166675533008643
$ws.Range("E16").Value = 0.6530552318172482
$ws.Range("F16").Value = 0.663354424801156
$ws.Range("G16").Value = 0.5737381753344266
$ws.Range("H16").Value = 0.4765477940424887
$ws.Range("I16").Value = 0.4896857341569927
$ws.Range("O16").Value = 2.069284775670098
$ws.Range("C17").Value = 0.009007556100804948
$ws.Range("D17").Value = 0.001131525787663179
$ws.Range("E17").Value = 0.6261350971456068
$ws.Range("F17").Value = 0.6449146304703959
$ws.Range("G17").Value = 0.5542159257400954
$ws.Range("H17").Value = 0.4689367221727991
$ws.Range("I17").Value = 0.4759578584805837
$ws.Range("O17").Value = 2.011823208150531
$ws.Range("C18").Value = 0.008824839663560624
$ws.Range("D18").Value = 0.001111337560338654
$ws.Range("E18").Value = 0.6106704135292489
$ws.Range("F18").Value = 0.6343404807311117
$ws.Range("G18").Value = 0.5430101100603792
$ws.Range("H18").Value = 0.4645848864146274
$ws.Range("I18").Value = 0.4680855176505787
$ws.Range("O18").Value = 1.978874020669139
$ws.Range("C19").Value = 0.008763019797221716
$ws.Range("D19").Value = 0.001104507180555103
$ws.Range("E19").Value = 0.6054375336137952
$ws.Range("F19").Value = 0.6307657236869346
$ws.Range("G19").Value = 0.5392199201961887
$ws.Range("H19").Value = 0.4631158635002066
$ws.Range("I19").Value = 0.4654241140248274
$ws.Range("O19").Value = 1.967735314281811
$ws.Range("C20").Value = 0.00904139404655524
$ws.Range("D20").Value = 0.001135264546594072
$ws.Range("E20").Value = 0.6289987995033499
$ws.Range("F20").Value = 0.6468742693821667
$ws.Range("G20").Value = 0.55629173112834
$ws.Range("H20").Value = 0.469744256791671
$ws.Range("I20").Value = 0.4774167718913844
$ws.Range("O20").Value = 2.017929604277867
$ws.Range("C21").Value = 0.009979788284667279
$ws.Range("D21").Value = 0.001238949162306824
$ws.Range("E21").Value = 0.7083950486137951
$ws.Range("F21").Value = 0.7013794881263919
$ws.Range("G21").Value = 0.6139253957688595
$ws.Range("H21").Value = 0.4923240747731654
$ws.Range("I21").Value = 0.5179929996918275
$ws.Range("O21").Value = 2.187788699467148
$ws.Range("C22").Value = 0.01059536607061062
$ws.Range("D22").Value = 0.001306964959901791
$ws.Range("E22").Value = 0.7604725738272293
$ws.Range("F22").Value = 0.737289768870852
$ws.Range("G22").Value = 0.6518005599495496
$ws.Range("H22").Value = 0.5073128822142507
$ws.Range("I22").Value = 0.5447246176007496
$ws.Range("O22").Value = 2.299715244977108
$ws.Range("C23").Value = 0.01026661325928302
$ws.Range("D23").Value = 0.001270640865492112
$ws.Range("E23").Value = 0.7326597997577977
$ws.Range("F23").Value = 0.7180972507515833
$ws.Range("G23").Value = 0.631566594094096
$ws.Range("H23").Value = 0.499291815127151
$ws.Range("I23").Value = 0.5304378499823486
$ws.Range("O23").Value = 2.23989370608416
$ws.Range("C24").Value = 0.009026095359509156
$ws.Range("D24").Value = 0.001133574192012432
$ws.Range("E24").Value = 0.6277040823208324
$ws.Range("F24").Value = 0.6459882321519217
$ws.Range("G24").Value = 0.5553532041648452
$ws.Range("H24").Value = 0.4693790960480442
$ws.Range("I24").Value = 0.476757134878909
$ws.Range("O24").Value = 2.015168634122176
$ws.Range("C25").Value = 0.007697946296993052
$ws.Range("D25").Value = 0.0009868378836532798
$ws.Range("E25").Value = 0.5152106233909137
$ws.Range("F25").Value = 0.5694254048093796
$ws.Range("G25").Value = 0.4381008095212735
$ws.Range("H25").Value = 0.4381008095212735
$ws.Range("I25").Value = 0.4197530610044709
$ws.Range("O25").Value = 1.776626452146161
